$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 75 (existing rows 75-83 shift down to 80-88)
$ws.Rows("75:79").Insert()

# Row 75
$ws.Cells.Item(75, 1).Value = 1
$ws.Cells.Item(75, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(75, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(75, 4).Value = 45218
$ws.Cells.Item(75, 5).Value = 15
$ws.Cells.Item(75, 6).Value = 'Fruta'
$ws.Cells.Item(75, 7).Value = 100104
$ws.Cells.Item(75, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(75, 9).Value = 100104005
$ws.Cells.Item(75, 10).Value = 'Pera'
$ws.Cells.Item(75, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(75, 12).Value = 'Primera'
$ws.Cells.Item(75, 13).Value = 250
$ws.Cells.Item(75, 14).Value = 19000
$ws.Cells.Item(75, 15).Value = 20000
$ws.Cells.Item(75, 16).Value = 19600
$ws.Cells.Item(75, 17).Value = '$/caja 20 kilos granel'
$ws.Cells.Item(75, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(75, 19).Value = 980
$ws.Cells.Item(75, 20).Value = 20

# Row 76
$ws.Cells.Item(76, 1).Value = 1
$ws.Cells.Item(76, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(76, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(76, 4).Value = 45218
$ws.Cells.Item(76, 5).Value = 15
$ws.Cells.Item(76, 6).Value = 'Fruta'
$ws.Cells.Item(76, 7).Value = 100104
$ws.Cells.Item(76, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(76, 9).Value = 100104005
$ws.Cells.Item(76, 10).Value = 'Pera'
$ws.Cells.Item(76, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(76, 12).Value = 'Segunda'
$ws.Cells.Item(76, 13).Value = 250
$ws.Cells.Item(76, 14).Value = 24000
$ws.Cells.Item(76, 15).Value = 25000
$ws.Cells.Item(76, 16).Value = 24500
$ws.Cells.Item(76, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(76, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(76, 19).Value = 1361
$ws.Cells.Item(76, 20).Value = 18

# Row 77
$ws.Cells.Item(77, 1).Value = 1
$ws.Cells.Item(77, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(77, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(77, 4).Value = 45218
$ws.Cells.Item(77, 5).Value = 15
$ws.Cells.Item(77, 6).Value = 'Fruta'
$ws.Cells.Item(77, 7).Value = 100104
$ws.Cells.Item(77, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(77, 9).Value = 100104005
$ws.Cells.Item(77, 10).Value = 'Pera'
$ws.Cells.Item(77, 11).Value = 'Packham''s Triumph'
$ws.Cells.Item(77, 12).Value = 'Tercera'
$ws.Cells.Item(77, 13).Value = 300
$ws.Cells.Item(77, 14).Value = 16000
$ws.Cells.Item(77, 15).Value = 17000
$ws.Cells.Item(77, 16).Value = 16500
$ws.Cells.Item(77, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(77, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(77, 19).Value = 917
$ws.Cells.Item(77, 20).Value = 18

# Row 78
$ws.Cells.Item(78, 1).Value = 1
$ws.Cells.Item(78, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(78, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(78, 4).Value = 45218
$ws.Cells.Item(78, 5).Value = 15
$ws.Cells.Item(78, 6).Value = 'Fruta'
$ws.Cells.Item(78, 7).Value = 100104
$ws.Cells.Item(78, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(78, 9).Value = 100104005
$ws.Cells.Item(78, 10).Value = 'Pera'
$ws.Cells.Item(78, 11).Value = 'Winter Nelis'
$ws.Cells.Item(78, 12).Value = 'Segunda'
$ws.Cells.Item(78, 13).Value = 300
$ws.Cells.Item(78, 14).Value = 24000
$ws.Cells.Item(78, 15).Value = 25000
$ws.Cells.Item(78, 16).Value = 24500
$ws.Cells.Item(78, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(78, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(78, 19).Value = 1361
$ws.Cells.Item(78, 20).Value = 18

# Row 79
$ws.Cells.Item(79, 1).Value = 1
$ws.Cells.Item(79, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(79, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(79, 4).Value = 45218
$ws.Cells.Item(79, 5).Value = 15
$ws.Cells.Item(79, 6).Value = 'Fruta'
$ws.Cells.Item(79, 7).Value = 100104
$ws.Cells.Item(79, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(79, 9).Value = 100104005
$ws.Cells.Item(79, 10).Value = 'Pera'
$ws.Cells.Item(79, 11).Value = 'Winter Nelis'
$ws.Cells.Item(79, 12).Value = 'Tercera'
$ws.Cells.Item(79, 13).Value = 300
$ws.Cells.Item(79, 14).Value = 16000
$ws.Cells.Item(79, 15).Value = 17000
$ws.Cells.Item(79, 16).Value = 16500
$ws.Cells.Item(79, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item(79, 18).Value = 'Región de O''Higgins'
$ws.Cells.Item(79, 19).Value = 917
$ws.Cells.Item(79, 20).Value = 18
